$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is used on the Overview sheet (E2:F3) and on each
# per-language sheet's "Status" column (C2:C3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# --- 2. Narrow the status columns ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de sheets: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
